$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking / percent-looking strings,
# then restore the default (Normal) style so no new cell formatting
# is introduced - matches the source workbook where these cells carry
# no explicit style.
$cells = @{
    'D2' = '330.43'
    'E2' = '7.20%'
    'D3' = '40.12'
    'E3' = '7.71%'
    'D4' = '5.402'
    'E4' = '5.34%'
    'D5' = '0.08107'
    'E5' = '3.54%'
    'D6' = '4.531'
    'E6' = '3.06%'
    'D7' = '8.658'
    'D8' = '1.925'
    'E8' = '2.35%'
    'D10' = '0.9464'
    'E10' = '2.74%'
    'D11' = '0.1361'
    'E11' = '25.91%'
    'D12' = '0.1974'
    'E12' = '4.20%'
    'D13' = '0.09322'
    'E13' = '5.09%'
    'D14' = '0.03566'
    'E14' = '7.70%'
    'D15' = '0.09589'
    'E15' = '-0.06%'
    'D16' = '0.001330'
    'E16' = '-3.43%'
    'E17' = '11.52%'
    'D18' = '3.361'
    'E18' = '-1.33%'
    'D19' = '0.3524'
    'E19' = '2.88%'
    'D20' = '7.236'
    'E20' = '14.99%'
    'E21' = '3.47%'
    'D23' = '0.04427'
    'E23' = '1.54%'
    'D24' = '0.001221'
    'E24' = '2.19%'
    'D25' = '0.004286'
    'E25' = '0.49%'
    'D26' = '0.0001200'
    'E26' = '-14.24%'
    'D27' = '0.0003990'
    'E27' = '-0.04%'
    'D39' = '0.02486'
    'E39' = '14.47%'
    'D40' = '0.05241'
    'E40' = '3.96%'
    'D41' = '0.007604'
    'E41' = '0.74%'
    'D42' = '0.1429'
    'E42' = '5.82%'
    'D43' = '0.009192'
    'E43' = '6.03%'
    'E44' = '9.61%'
    'D45' = '0.01090'
    'E45' = '38.10%'
    'D46' = '0.00006598'
    'E46' = '1.00%'
    'E47' = '0.10%'
    'D48' = '0.002400'
    'E48' = '139.51%'
    'E49' = '1.53%'
    'D50' = '0.00002100'
    'E50' = '0.10%'
    'D51' = '0.0002000'
    'E51' = '0.10%'
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
    $rng.Style = "Normal"
}
